$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")
$ws.Activate()

# ---------------------------------------------------------------------------
# 1) View state: scroll the frozen pane back to the top (topLeftCell B3) and
#    move the active selection in the bottom-right pane to A2.
# ---------------------------------------------------------------------------
[void]$ws.Range("B3").Select()
$excel.ActiveWindow.FreezePanes = $true
[void]$ws.Range("A2").Select()

# ---------------------------------------------------------------------------
# 2) Data updates for rows 217-225 (only the raw input cells change; columns
#    B, H, J and K are formulas and recompute automatically).
# ---------------------------------------------------------------------------
$ws.Range("P217").Value2 = 581
$ws.Range("P218").Value2 = 579

$ws.Range("N219").Value2 = 136
$ws.Range("P219").Value2 = 615

$ws.Range("N220").Value2 = 137
$ws.Range("P220").Value2 = 688

$ws.Range("N221").Value2 = 138
$ws.Range("P221").Value2 = 774

$ws.Range("N222").Value2 = 144
$ws.Range("O222").Value2 = 463
$ws.Range("P222").Value2 = 855

$ws.Range("C223").Value2 = 34
$ws.Range("N223").Value2 = 167
$ws.Range("O223").Value2 = 433
$ws.Range("P223").Value2 = 907

$ws.Range("C224").Value2 = 60
$ws.Range("N224").Value2 = 213
$ws.Range("O224").Value2 = 480
$ws.Range("P224").Value2 = 867

$ws.Range("C225").Value2 = 58
$ws.Range("D225").Value2 = 1
$ws.Range("G225").Value2 = 9
$ws.Range("N225").Value2 = 260
$ws.Range("O225").Value2 = 542
$ws.Range("P225").Value2 = 818

# ---------------------------------------------------------------------------
# 3) Row 226 was a placeholder (blank inputs, formulas returning "") and is
#    now fully populated with real data, matching the formatting already
#    used by row 225 (L/M switch from style 18/19 to style 7).
# ---------------------------------------------------------------------------
$ws.Range("L225:M225").Copy()
$ws.Range("L226:M226").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("C226").Value2 = 0
$ws.Range("D226").Value2 = 0
$ws.Range("E226").Value2 = 1
$ws.Range("F226").Value2 = 1
$ws.Range("G226").Value2 = 9
$ws.Range("I226").Value2 = 0
$ws.Range("L226").Value2 = 0
$ws.Range("M226").Value2 = 0
$ws.Range("N226").Value2 = 244
$ws.Range("O226").Value2 = 521
$ws.Range("P226").Value2 = 724
